# DEBUG: push exchange event into timeline
# Extend the basketball timeline sheet with six overtime period blocks
# (延長賽一 .. 延長賽六), each mirroring the existing quarter layout:
#   [Period, Player, 進攻方式, 出手方式, 結果, 罰球, Pts] + 1 blank gap column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$subHeaders = @("Player", "進攻方式", "出手方式", "結果", "罰球", "Pts")
$otNames = @("延長賽一", "延長賽二", "延長賽三", "延長賽四", "延長賽五", "延長賽六")

# First existing block starts at column A (1); each block is 7 columns wide
# followed by 1 blank column (8 total). Existing quarters occupy blocks 0-3
# (A.. through Y..), so the new overtime blocks continue on immediately
# after, starting at block index 4 (column AG).
$startBlock = 4

for ($i = 0; $i -lt $otNames.Length; $i++) {
    $blockIndex = $startBlock + $i
    $headerCol = ($blockIndex * 8) + 1

    $ws.Cells.Item(1, $headerCol).Value = $otNames[$i]

    for ($j = 0; $j -lt $subHeaders.Length; $j++) {
        $ws.Cells.Item(1, $headerCol + 1 + $j).Value = $subHeaders[$j]
    }
}

# Mirror the author's view-state: scrolled over to the new columns with
# AH3 as the active cell.
$ws.Range("AH3").Select()
